# Auto-generated edit script applying market price / profit updates
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 7500
$ws.Range("I43").Value = 7500
$ws.Range("K43").Value = 7500
$ws.Range("M43").Value = -7431

$ws.Range("H112").Value = 2173.647
$ws.Range("J112").Value = 2173.647
$ws.Range("L112").Value = 6520.941
$ws.Range("N112").Value = -8736.940999999999

$ws.Range("H132").Value = 1752.0303
$ws.Range("I132").Value = 1793.9
$ws.Range("K132").Value = 5381.700000000001
$ws.Range("M132").Value = -2851.700000000001

$ws.Range("H135").Value = 1079.3429
$ws.Range("J135").Value = 1520.625
$ws.Range("L135").Value = 13685.625
$ws.Range("N135").Value = -18755.625

$ws.Range("H137").Value = 2453.2083
$ws.Range("J137").Value = 2929.6667
$ws.Range("L137").Value = 8789.000100000001
$ws.Range("N137").Value = -13889.0001

$ws.Range("H138").Value = 5771
$ws.Range("I138").Value = 6592.778
$ws.Range("K138").Value = 19778.334
$ws.Range("M138").Value = -14638.334

$ws.Range("H141").Value = 3904.5
$ws.Range("I141").Value = 3262.8572
$ws.Range("J141").Value = 5401.6665
$ws.Range("K141").Value = 9788.571599999999
$ws.Range("L141").Value = 16204.9995
$ws.Range("M141").Value = -4608.571599999999
$ws.Range("N141").Value = -26564.9995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1789
$ws.Range("J2").Value = 2674.3333
$ws.Range("L2").Value = 2674.3333
$ws.Range("N2").Value = -2900.3333

$ws.Range("H32").Value = 5119.018
$ws.Range("I32").Value = 2883.068
$ws.Range("J32").Value = 14062.818
$ws.Range("K32").Value = 2883.068
$ws.Range("L32").Value = 14062.818
$ws.Range("M32").Value = -2596.068
$ws.Range("N32").Value = -14636.818

$ws.Range("H116").Value = 1789
$ws.Range("J116").Value = 2674.3333
$ws.Range("L116").Value = 2674.3333
$ws.Range("N116").Value = -7262.3333

$ws.Range("H132").Value = 1481.0714
$ws.Range("I132").Value = 1373.909
$ws.Range("J132").Value = 1874
$ws.Range("K132").Value = 4121.727000000001
$ws.Range("L132").Value = 5622
$ws.Range("M132").Value = -1591.727000000001
$ws.Range("N132").Value = -10682

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1789
$ws.Range("J3").Value = 2674.3333
$ws.Range("L3").Value = 2674.3333
$ws.Range("N3").Value = -2902.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 681
$ws.Range("I107").Value = 621.5
$ws.Range("K107").Value = 621.5
$ws.Range("M107").Value = 1298.5

$ws.Range("H132").Value = 2722.9583
$ws.Range("I132").Value = 2406.6086
$ws.Range("K132").Value = 7219.825800000001
$ws.Range("M132").Value = -4689.825800000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 71453.07000000001
$ws.Range("I2").Value = 125017.125
$ws.Range("K2").Value = 750102.75
$ws.Range("M2").Value = -749989.75

$ws.Range("H17").Value = 910.7143
$ws.Range("I17").Value = 90.5
$ws.Range("J17").Value = 1238.8
$ws.Range("K17").Value = 271.5
$ws.Range("L17").Value = 3716.4
$ws.Range("M17").Value = -102.5
$ws.Range("N17").Value = -4054.4

$ws.Range("H32").Value = 2401
$ws.Range("I32").Value = 802
$ws.Range("J32").Value = 4000
$ws.Range("K32").Value = 2406
$ws.Range("L32").Value = 12000
$ws.Range("M32").Value = -2123
$ws.Range("N32").Value = -12566

$ws.Range("H34").Value = 2164.889
$ws.Range("I34").Value = 871
$ws.Range("J34").Value = 3200
$ws.Range("K34").Value = 2613
$ws.Range("L34").Value = 9600
$ws.Range("M34").Value = -2529
$ws.Range("N34").Value = -9768

$ws.Range("H55").Value = 1000000
$ws.Range("I55").Value = 1000000
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 3000000
$ws.Range("L55").Value = 0
$ws.Range("M55").Value = -2999823
$ws.Range("N55").ClearContents()

$ws.Range("H128").Value = 3979890
$ws.Range("I128").Value = 3979890
$ws.Range("K128").Value = 11939670
$ws.Range("M128").Value = -11934690

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 6792.2856
$ws.Range("J80").Value = 9749
$ws.Range("L80").Value = 9749
$ws.Range("N80").Value = -11745

$ws.Range("H83").Value = 6792.2856
$ws.Range("J83").Value = 9749
$ws.Range("L83").Value = 48745
$ws.Range("N83").Value = -58729

$ws.Range("H126").Value = 4187.6665
$ws.Range("I126").Value = 3037.3333
$ws.Range("K126").Value = 9111.999899999999
$ws.Range("M126").Value = -6641.999899999999

$ws.Range("H132").Value = 2512
$ws.Range("I132").Value = 1711.9231
$ws.Range("J132").Value = 3457.5454
$ws.Range("K132").Value = 5135.7693
$ws.Range("L132").Value = 10372.6362
$ws.Range("M132").Value = -2605.7693
$ws.Range("N132").Value = -15432.6362

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3356.2
$ws.Range("I40").Value = 2945.25
$ws.Range("K40").Value = 2945.25
$ws.Range("M40").Value = -2809.25

$ws.Range("H46").Value = 3573.25
$ws.Range("I46").Value = 2296
$ws.Range("J46").Value = 4339.6
$ws.Range("K46").Value = 2296
$ws.Range("L46").Value = 4339.6
$ws.Range("M46").Value = -2108
$ws.Range("N46").Value = -4715.6

$ws.Range("H74").Value = 50000
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents()

$ws.Range("H77").Value = 50000
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents()

$ws.Range("H122").Value = 5560.7
$ws.Range("I122").Value = 5997.5454
$ws.Range("J122").Value = 5026.778
$ws.Range("K122").Value = 17992.6362
$ws.Range("L122").Value = 15080.334
$ws.Range("M122").Value = -15542.6362
$ws.Range("N122").Value = -19980.334

$ws.Range("H132").Value = 3361.1714
$ws.Range("I132").Value = 2576.875
$ws.Range("K132").Value = 7730.625
$ws.Range("M132").Value = -5200.625

$ws.Range("H133").Value = 105000
$ws.Range("J133").Value = 105000
$ws.Range("L133").Value = 105000
$ws.Range("N133").Value = -110060

$ws.Range("H136").Value = 10000
$ws.Range("J136").Value = 10000
$ws.Range("L136").Value = 30000
$ws.Range("N136").Value = -35100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H138").Value = 190000
$ws.Range("J138").Value = 190000
$ws.Range("L138").Value = 190000
$ws.Range("N138").Value = -200280
